$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AND OR")
$ws1.Range("P6").Value = "test"
Write-Host "P6:" $ws1.Range("P6").Value()
$ws1.Range("P6").ClearContents()
Write-Host "P6 after clear:" $ws1.Range("P6").Value()
Write-Host "P6 style:" $ws1.Range("P6").Style
$ws1.Range("P6").Style = $ws1.Range("Q6").Style
